{"js": "// Fix the \"None\" paragraph under \"Prerequisite tasks\": it was incorrectly\n// authored using the Heading1 style (with an explicit 22-half-point / 11pt\n// run size override) instead of plain body text, and it carried a stray\n// bookmark. Convert it to a normal paragraph and drop the bookmark \u2014 the\n// remaining bookmarks (Support Information / Detailed Instructions /\n// Document Management) are then renumbered down by one automatically.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph whose entire text is exactly \"None\" (the\n// Prerequisite-tasks placeholder, not any other occurrence of the word).\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"None\") {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  const range = targetParagraph.getRange();\n\n  // Replace the paragraph's OOXML with a plain-paragraph version: no\n  // pStyle (=> Normal style), no sz/szCs override, and no bookmark.\n  const replacementPackage =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:pPr><w:rPr/></w:pPr>' +\n    '<w:r><w:rPr><w:rtl w:val=\"0\"/></w:rPr><w:t xml:space=\"preserve\">None</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n\n  range.insertOoxml(replacementPackage, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fix the \"None\" paragraph under \"Prerequisite tasks\": it was incorrectly\n# authored using the Heading1 style (with an explicit 22-half-point / 11pt\n# run size override) instead of plain body text, and it carried a stray\n# bookmark. Convert it to a normal paragraph and drop the bookmark \u2014 the\n# remaining bookmarks (Support Information / Detailed Instructions /\n# Document Management) are then renumbered down by one automatically.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"None\"\n$rng.Find.MatchWholeWord = $true\n$rng.Find.MatchCase = $true\n$found = $rng.Find.Execute()\n\nif ($found) {\n    # Resolve the paragraph that owns the found text. Operating on the\n    # paragraph's own Range (rather than the collapsed Find range) makes\n    # InsertXML replace the paragraph's content in place instead of just\n    # inserting a sibling run.\n    $target = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {\n            $target = $p\n            break\n        }\n    }\n\n    if ($target -ne $null) {\n        $r = $target.Range\n\n        # Replace the paragraph's OOXML with a plain-paragraph version: no\n        # pStyle (=> Normal style), no sz/szCs override, and no bookmark.\n        $pkgXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n            '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n            '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n            '<pkg:xmlData>' +\n            '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' +\n            '<w:p>' +\n            '<w:pPr><w:rPr/></w:pPr>' +\n            '<w:r><w:rPr><w:rtl w:val=\"0\"/></w:rPr><w:t xml:space=\"preserve\">None</w:t></w:r>' +\n            '</w:p>' +\n            '</w:body>' +\n            '</w:document>' +\n            '</pkg:xmlData>' +\n            '</pkg:part>' +\n            '</pkg:package>'\n\n        $r.InsertXML($pkgXml)\n    }\n}\n"}
